# "Added Paths as Dropdown, For making the process of choosing the paths easier!"
# Append a new data row (row 3) that mirrors row 2, then edit the
# Time Period / Comments / Win-Lost values for the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 into row 3 first (via Copy) so the Date cell keeps its
# original text type instead of being re-parsed/auto-converted to a date
# serial number when retyped.
$ws.Range("A2:D2").Copy($ws.Range("A3:D3"))

# Now overwrite the new row's Time Period, Comments and Win/Lost entries.
$ws.Range("B3").Value = "2-3"
$ws.Range("C3").Value = "hitdujgrfjdthgfdsfgf"
$ws.Range("D3").Value = "Win"
